$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F ("dSF") values for the listed rows, per repulled data.
$ws.Range("F7").Value = -4
$ws.Range("F14").Value = -3
$ws.Range("F15").Value = 1
$ws.Range("F17").Value = 6
$ws.Range("F20").Value = -5
$ws.Range("F22").Value = -3
$ws.Range("F23").Value = 1
$ws.Range("F24").Value = 1
$ws.Range("F32").Value = 0
$ws.Range("F33").Value = 7
$ws.Range("F34").Value = 3
$ws.Range("F35").Value = -3
$ws.Range("F36").Value = -6
$ws.Range("F38").Value = 0
$ws.Range("F42").Value = -1
